# Updates the Price (D) and Volume(1h) (E) columns for the cryptos
# report rows, matching the refreshed source data.
#
# D-column values are prefixed with a leading apostrophe so Excel
# keeps them as literal text (e.g. preserves "1.010" / "11.00"
# instead of normalising to the number 1.01 / 11) -- matching the
# original inline-string cell content exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.648.98"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").Value = "'1.844.32"
$ws.Range("E3").Value = "  -0.15%  "

$ws.Range("D4").Value = "'1.010"
$ws.Range("E4").Value = "  -1.95%  "

$ws.Range("D5").Value = "'317.17"
$ws.Range("E5").Value = "  -1.34%  "

$ws.Range("D6").Value = "'1.008"
$ws.Range("E6").Value = "  -1.78%  "

$ws.Range("D7").Value = "'0.4293"
$ws.Range("E7").Value = "  -1.84%  "

$ws.Range("D8").Value = "'0.3737"
$ws.Range("E8").Value = "  -1.20%  "

$ws.Range("D9").Value = "'0.07310"
$ws.Range("E9").Value = "  -0.79%  "

$ws.Range("D10").Value = "'0.8721"
$ws.Range("E10").Value = "  -0.92%  "

$ws.Range("D11").Value = "'21.43"
$ws.Range("E11").Value = "  -0.31%  "

$ws.Range("D12").Value = "'1.847.50"
$ws.Range("E12").Value = "  -0.40%  "

$ws.Range("D13").Value = "'6.710"
$ws.Range("E13").Value = "  +0.19%  "

$ws.Range("D14").Value = "'5.403"
$ws.Range("E14").Value = "  -1.66%  "

$ws.Range("D15").Value = "'0.07122"
$ws.Range("E15").Value = "  -0.10%  "

$ws.Range("D16").Value = "'88.79"
$ws.Range("E16").Value = "  +4.45%  "

$ws.Range("D17").Value = "'1.010"
$ws.Range("E17").Value = "  -2.13%  "

$ws.Range("D18").Value = "'0.000008979"
$ws.Range("E18").Value = "  -0.66%  "

$ws.Range("E19").Value = "  -1.68%  "

$ws.Range("D20").Value = "'15.38"
$ws.Range("E20").Value = "  -0.26%  "

$ws.Range("D21").Value = "'27.663.01"
$ws.Range("E21").Value = "  +0.14%  "

$ws.Range("D22").Value = "'5.193"

$ws.Range("D23").Value = "'11.00"
$ws.Range("E23").Value = "  -2.55%  "

$ws.Range("D24").Value = "'2.070.55"
$ws.Range("E24").Value = "  -0.40%  "

$ws.Range("D25").Value = "'1.967"
$ws.Range("E25").Value = "  -3.32%  "

$ws.Range("D26").Value = "'154.70"
$ws.Range("E26").Value = "  -1.83%  "

$ws.Range("D27").Value = "'18.55"
$ws.Range("E27").Value = "  -0.64%  "

$ws.Range("D28").Value = "'2.159"
$ws.Range("E28").Value = "  +7.96%  "

$ws.Range("D29").Value = "'5.328"
$ws.Range("E29").Value = "  +0.18%  "

$ws.Range("D30").Value = "'117.72"
$ws.Range("E30").Value = "  +0.07%  "

$ws.Range("D31").Value = "'0.08904"
$ws.Range("E31").Value = "  -1.25%  "

$ws.Range("D32").Value = "'1.218"
$ws.Range("E32").Value = "  +1.11%  "

$ws.Range("D33").Value = "'0.7737"
$ws.Range("E33").Value = "  +0.64%  "

$ws.Range("D34").Value = "'4.524"
$ws.Range("E34").Value = "  -0.32%  "

$ws.Range("D35").Value = "'2.899"
$ws.Range("E35").Value = "  -2.97%  "

$ws.Range("E36").Value = "  -1.72%  "

$ws.Range("E37").Value = "  -1.18%  "

$ws.Range("D38").Value = "'0.01973"
$ws.Range("E38").Value = "  +0.25%  "

$ws.Range("D39").Value = "'0.05299"
$ws.Range("E39").Value = "  +0.75%  "

$ws.Range("D40").Value = "'2.885"
$ws.Range("E40").Value = "  +1.68%  "

$ws.Range("D41").Value = "'7.153"
$ws.Range("E41").Value = "  +4.85%  "

$ws.Range("D42").Value = "'0.1687"
$ws.Range("E42").Value = "  +1.25%  "

$ws.Range("E43").Value = "  -0.98%  "

$ws.Range("D44").Value = "'8.755"
$ws.Range("E44").Value = "  +0.11%  "

$ws.Range("D45").Value = "'10.64"
$ws.Range("E45").Value = "  -0.67%  "

$ws.Range("E46").Value = "  -2.85%  "

$ws.Range("D47").Value = "'0.4740"
$ws.Range("E47").Value = "  +1.20%  "

$ws.Range("D48").Value = "'0.06450"
$ws.Range("E48").Value = "  -2.41%  "

$ws.Range("D49").Value = "'1.009"
$ws.Range("E49").Value = "  -1.88%  "

$ws.Range("D50").Value = "'1.683"
$ws.Range("E50").Value = "  -0.85%  "

$ws.Range("D51").Value = "'1.838"
$ws.Range("E51").Value = "  -2.41%  "

